# Update a batch of imputed numeric results in the RandomForest result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E6").Value = 16.38680000000002
$ws.Range("D7").Value = -7.166100000000005
$ws.Range("B9").Value = 6.311499999999998
$ws.Range("D12").Value = -7.037400000000003
$ws.Range("D14").Value = -7.728900000000005
$ws.Range("E15").Value = 16.13190000000002
$ws.Range("B18").Value = 7.2094
$ws.Range("B20").Value = 8.991199999999997
$ws.Range("D26").Value = -8.388500000000004
$ws.Range("B27").Value = 5.836400000000005
$ws.Range("D27").Value = -8.774700000000001
$ws.Range("D29").Value = -7.240899999999999
$ws.Range("E33").Value = 17.31890000000002
$ws.Range("B35").Value = 8.715
$ws.Range("E35").Value = 16.44190000000001
$ws.Range("D37").Value = -7.7385
$ws.Range("D38").Value = -7.621199999999996
$ws.Range("E38").Value = 16.65129999999999
$ws.Range("E43").Value = 17.23160000000001
$ws.Range("E44").Value = 16.06299999999999
$ws.Range("E47").Value = 16.032
$ws.Range("D51").Value = -8.594299999999999
$ws.Range("E51").Value = 16.3626
$ws.Range("D52").Value = -7.249699999999998
$ws.Range("D55").Value = -8.627799999999999
$ws.Range("E57").Value = 16.0212
$ws.Range("E63").Value = 18.11780000000001
$ws.Range("B69").Value = 5.373599999999995
$ws.Range("D69").Value = -7.126999999999996
$ws.Range("D70").Value = -8.350200000000003
$ws.Range("E70").Value = 16.9766
$ws.Range("B76").Value = 5.199299999999998
$ws.Range("B78").Value = 9.789399999999999
$ws.Range("D81").Value = -7.880200000000007
$ws.Range("B82").Value = 6.461599999999999
$ws.Range("B83").Value = 5.577900000000001
$ws.Range("D83").Value = -9.059299999999997
$ws.Range("E88").Value = 16.4273
$ws.Range("B93").Value = 6.7385
$ws.Range("E99").Value = 16.58459999999999
$ws.Range("D102").Value = -7.879599999999998
